$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Gold color used for currency symbol "đ" (FFBB934D) as a COM color long (BGR-ish packed value used by .Font.Color)
$goldColor = 5084091   # 0x4D93BB -> R=0xBB,G=0x93,B=0x4D
$blackColor = 0

function Set-PriceCell($cell, $amountText) {
    $cell.Value = $amountText + [char]0x00A0 + [char]0x0111
    $len = $amountText.Length
    $sp = $cell.Characters($len + 1, 1)
    $sp.Font.Size = 8
    $sp.Font.Color = $blackColor
    $sp.Font.Name = "Arial"
    $dong = $cell.Characters($len + 2, 1)
    $dong.Font.Size = 17
    $dong.Font.Color = $goldColor
    $dong.Font.Name = "Arial"
}

$rows = @(
    @{ Row=82; Height=27.6; B="GD0000W000105";  C="Dây chuyền Vàng trắng 10K PNJ dây đan kiểu chữ cong 0000W000105";       Price="1.080.000"; F="https://www.pnj.com.vn/day-chuyen-pnj-vang-trang-10k-8094.html";                 CStyle6=$false },
    @{ Row=83; Height=21;   B="GD0000W000220";  C="Dây chuyền Vàng trắng 10K PNJ 0000W000220";                            Price="1.322.000"; F="https://www.pnj.com.vn/day-chuyen-pnj-vang-trang-10k-gdmrwkxx025.006.html";      CStyle6=$false },
    @{ Row=84; Height=21;   B=" GD0000W000177"; C="Dây chuyền Vàng trắng 10K PNJ 0000W000177";                            Price="1.717.000"; F="https://www.pnj.com.vn/day-chuyen-vang-trang-10k-6813.html";                     CStyle6=$false },
    @{ Row=85; Height=21;   B="GD0000Y000711";  C="Dây chuyền Vàng 18K PNJ kiểu dây đan dập chữ S xoắn suốt 0000Y000711"; Price="2.341.000"; F="https://www.pnj.com.vn/day-chuyen-vang-18k-6555.html?";                          CStyle6=$true  },
    @{ Row=86; Height=21;   B="GD0000W000277";  C="Dây chuyền Vàng trắng Ý 18K PNJ 0000W000277";                          Price="2.175.000"; F="https://www.pnj.com.vn/day-chuyen-pnj-vang-trang-y-18k-8727.html";              CStyle6=$true  },
    @{ Row=87; Height=21;   B="GD0000Y000256";  C="Dây chuyền Vàng 18K PNJ dây đan kiểu chữ cong 0000Y000256";            Price="1.862.000"; F="https://www.pnj.com.vn/day-chuyen-pnj-vang-18k-day-dan-kieu-chu-cong-vi.html";   CStyle6=$true  },
    @{ Row=88; Height=21;   B="GCXMXMY000021";  C="Dây cổ Vàng 18K đính đá CZ PNJ XMXMY000021";                          Price="6.713.000"; F="https://www.pnj.com.vn/day-co-pnj-vang-18k-dinh-da-cz-7256.html";               CStyle6=$false },
    @{ Row=89; Height=21;   B="GCXMXMY000019";  C="Dây cổ Vàng 18K đính đá CZ PNJ XMXMY000019";                          Price="6.627.000"; F="https://www.pnj.com.vn/day-co-pnj-vang-18k-dinh-da-cz-ya59993.102.html";        CStyle6=$false }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    # Copy formatting for B:F from the row directly above the block (row 81), which already
    # carries the B=1 / C=2 / D=3 / F=(none) style combination used throughout the sheet.
    $ws.Range("B81:F81").Copy()
    $ws.Range("B$rowNum`:F$rowNum").PasteSpecial(-4122)
    $excel.CutCopyMode = 0

    if ($r.CStyle6) {
        # Rows whose product-name cell uses the non-wrapping style (style id 6) -- copy that
        # specific cell format from C4, which already uses it.
        $ws.Range("C4").Copy()
        $ws.Range("C$rowNum").PasteSpecial(-4122)
        $excel.CutCopyMode = 0
    }

    $ws.Rows.Item($rowNum).RowHeight = $r.Height

    $ws.Cells.Item($rowNum, 2).Value = $r.B
    $ws.Cells.Item($rowNum, 3).Value = $r.C
    $ws.Cells.Item($rowNum, 6).Value = $r.F

    Set-PriceCell $ws.Cells.Item($rowNum, 4) $r.Price
}

# The author's selection/scroll position changed after adding the rows above.
$excel.ActiveWindow.ScrollRow = 82
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("D92").Select() | Out-Null
